$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# "Recorded By" (column G) lists the users who recorded a session, as a
# comma-separated string. This processing pass swaps the order of the two
# names for every two-author cell EXCEPT the "System, backup@backdoor.com"
# pairing, which keeps its original order.
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Range("G$r")
    $val = $cell.Value2
    if ($null -eq $val) { continue }

    $allParts = $val -split ", "
    if ($allParts.Length -eq 2 -and $val -ne "System, backup@backdoor.com") {
        $cell.Value2 = "$($allParts[1]), $($allParts[0])"
    }
}
